$p = $ppt.ActivePresentation
$p.Slides.Item(36).Delete()
